# Applies odds updates for Jogos_do_Dia_Betfair_Back_Lay_2025-10-28.xlsx
# Only numeric odds values change; no structural changes to rows/columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("J2").Value = 5.1

# Row 3
$ws.Range("G3").Value = 8.199999999999999
$ws.Range("H3").Value = 1.55
$ws.Range("I3").Value = 1.56
$ws.Range("J3").Value = 4.3
$ws.Range("K3").Value = 4.4
$ws.Range("N3").Value = 3.6
$ws.Range("P3").Value = 1.83
$ws.Range("Q3").Value = 2.12
$ws.Range("Y3").Value = 7.8
$ws.Range("Z3").Value = 9.199999999999999
$ws.Range("AB3").Value = 26
$ws.Range("AC3").Value = 10

# Row 4
$ws.Range("F4").Value = 2.84
$ws.Range("G4").Value = 3.2
$ws.Range("H4").Value = 2.64
$ws.Range("I4").Value = 2.92
$ws.Range("J4").Value = 3.1

# Row 5
$ws.Range("Q5").Value = 2.32

# Row 6
$ws.Range("F6").Value = 2.26
$ws.Range("G6").Value = 2.46
$ws.Range("J6").Value = 3.1
$ws.Range("K6").Value = 3.6
$ws.Range("Q6").Value = 2.16

# Row 7
$ws.Range("F7").Value = 4.2
$ws.Range("G7").Value = 5.1
$ws.Range("H7").Value = 1.96
$ws.Range("I7").Value = 2.42
$ws.Range("J7").Value = 2.88
$ws.Range("K7").Value = 3.7
$ws.Range("P7").Value = 1.71
$ws.Range("Q7").Value = 2.18

# Row 8
$ws.Range("G8").Value = 2.2
$ws.Range("H8").Value = 4
$ws.Range("I8").Value = 4.6
$ws.Range("K8").Value = 3.7
$ws.Range("P8").Value = 1.74

# Row 10
$ws.Range("J10").Value = 4.7

# Row 11
$ws.Range("F11").Value = 2.22
$ws.Range("G11").Value = 2.66
$ws.Range("H11").Value = 2.96

# Row 12
$ws.Range("F12").Value = 2.74
$ws.Range("M12").Value = 1.07
$ws.Range("N12").Value = 4.2
$ws.Range("O12").Value = 1.3
$ws.Range("R12").Value = 1.39
$ws.Range("X12").Value = 15.5
$ws.Range("Y12").Value = 13
$ws.Range("Z12").Value = 1000
$ws.Range("AA12").Value = 50
$ws.Range("AB12").Value = 12.5
$ws.Range("AC12").Value = 7.8
$ws.Range("AD12").Value = 13.5
$ws.Range("AE12").Value = 1000
$ws.Range("AI12").Value = 48
$ws.Range("AJ12").Value = 50
$ws.Range("AL12").Value = 1000
$ws.Range("AO12").Value = 1000

# Row 13
$ws.Range("G13").Value = 1.97
$ws.Range("H13").Value = 3.45
$ws.Range("I13").Value = 7.4
$ws.Range("J13").Value = 3
$ws.Range("Q13").Value = 1.02

# Row 14
$ws.Range("P14").Value = 2.16
$ws.Range("Q14").Value = 1.62

# Row 15
$ws.Range("F15").Value = 1.91
$ws.Range("G15").Value = 2.14
$ws.Range("H15").Value = 3.45
$ws.Range("I15").Value = 6
$ws.Range("J15").Value = 2.9
$ws.Range("K15").Value = 4.1
$ws.Range("P15").Value = 1.5
$ws.Range("Q15").Value = 2.18
